$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.755.65'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.603.08'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.93'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.69'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.827.83'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.601.82'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.06'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '209.98'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.15'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.23'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -4.23%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '143.59'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.09'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.291.05'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.19'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +13.21%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.833'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.780'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '63.09'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.739.96'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0517'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.42'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.35%  '
